$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Donovan Mitchell",   "PG,SG",    "Cleveland Cavaliers"),
    @("Malik Beasley",      "SG,SF",    "Detroit Pistons"),
    @("De'Andre Hunter",    "SF,PF",    "Cleveland Cavaliers"),
    @("Jaden McDaniels",    "SF,PF",    "Minnesota Timberwolves"),
    @("Devin Vassell",      "SG,SF",    "San Antonio Spurs"),
    @("Josh Hart",          "SG,SF,PF", "New York Knicks"),
    @("Kristaps Porzingis", "PF,C",     "Boston Celtics"),
    @("Dyson Daniels",      "PG,SG,SF", "Atlanta Hawks"),
    @("Alperen Sengün",     "C",        "Houston Rockets"),
    @("Domantas Sabonis",   "C",        "Sacramento Kings"),
    @("Toumani Camara",     "SG,SF,PF", "Portland Trail Blazers"),
    @("Kelly Oubre Jr.",    "SG,SF",    "Philadelphia 76ers"),
    @("Michael Porter Jr.", "SF,PF",    "Denver Nuggets"),
    @("Aaron Wiggins",      "SG,SF",    "Oklahoma City Thunder"),
    @("Moussa Diabate",     "C",        "Charlotte Hornets"),
    @("Julius Randle",      "PF,C",     "Minnesota Timberwolves"),
    @("Cam Thomas",         "SG,SF",    "Brooklyn Nets"),
    @("Tari Eason",         "SF,PF",    "Houston Rockets")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
